$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" '25.848.25'
Set-TextValue "E2" '  -0.04%  '

Set-TextValue "D3" '1.642.14'
Set-TextValue "E3" '  +0.37%  '

Set-TextValue "E4" '  -0.21%  '

Set-TextValue "D5" '216.30'
Set-TextValue "E5" '  +0.37%  '

Set-TextValue "D6" '0.505'
Set-TextValue "E6" '  -0.58%  '

Set-TextValue "E8" '  -0.35%  '

Set-TextValue "D10" '19.77'
Set-TextValue "E10" '  -1.92%  '

Set-TextValue "E11" '  +1.69%  '

Set-TextValue "E12" '  +0.42%  '

Set-TextValue "D13" '1.868.05'
Set-TextValue "E13" '  +0.33%  '

Set-TextValue "D14" '1.641.15'
Set-TextValue "E14" '  -0.11%  '

Set-TextValue "E15" '  -0.12%  '

Set-TextValue "D16" '0.0₃0768'

Set-TextValue "D17" '63.13'
Set-TextValue "E17" '  -0.22%  '

Set-TextValue "D18" '25.882.48'
Set-TextValue "E18" '  +0.14%  '

Set-TextValue "E19" '  -0.20%  '

Set-TextValue "D20" '4.48'
Set-TextValue "E20" '  +2.28%  '

Set-TextValue "D21" '193.14'
Set-TextValue "E21" '  -0.55%  '

Set-TextValue "D22" '9.99'
Set-TextValue "E22" '  +0.71%  '

Set-TextValue "D23" '6.35'
Set-TextValue "E23" '  +2.52%  '

Set-TextValue "E24" '  +4.99%  '

Set-TextValue "E25" '  -0.18%  '

Set-TextValue "D26" '142.12'
Set-TextValue "E26" '  +2.62%  '

Set-TextValue "E27" '  +0.30%  '

Set-TextValue "D28" '6.97'
Set-TextValue "E28" '  +1.94%  '

Set-TextValue "D29" '15.57'
Set-TextValue "E29" '  +0.10%  '

Set-TextValue "E30" '  -0.06%  '

Set-TextValue "D31" '0.0496'
Set-TextValue "E31" '  +0.50%  '

Set-TextValue "E32" '  +1.09%  '

Set-TextValue "D33" '3.26'
Set-TextValue "E33" '  +0.52%  '

Set-TextValue "D34" '1.59'
Set-TextValue "E34" '  +0.84%  '

Set-TextValue "E35" '  -0.46%  '

Set-TextValue "D36" '0.911'
Set-TextValue "E36" '  +0.77%  '

Set-TextValue "D37" '1.134.75'
Set-TextValue "E37" '  +1.24%  '

Set-TextValue "B38" 'MXToken'
Set-TextValue "C38" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D38" '2.53'
Set-TextValue "E38" '  -1.99%  '

Set-TextValue "B39" 'ImmutableX'
Set-TextValue "C39" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D39" '0.548'
Set-TextValue "E39" '  -0.15%  '

Set-TextValue "E40" '  -0.51%  '

Set-TextValue "E41" '  +0.03%  '

Set-TextValue "E42" '  +1.32%  '

Set-TextValue "D43" '100.81'
Set-TextValue "E43" '  +1.35%  '

Set-TextValue "D44" '0.808'
Set-TextValue "E44" '  +0.85%  '

Set-TextValue "D45" '1.777.00'
Set-TextValue "E45" '  +0.29%  '

Set-TextValue "E46" '  +0.00%  '

Set-TextValue "E47" '  +0.07%  '

Set-TextValue "B48" 'RenderToken'
Set-TextValue "C48" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D48" '1.45'
Set-TextValue "E48" '  +6.43%  '

Set-TextValue "B49" 'Mantle'
Set-TextValue "C49" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D49" '0.419'
Set-TextValue "E49" '  -1.16%  '

Set-TextValue "E50" '  -0.06%  '

Set-TextValue "B51" 'EnergySwap'
Set-TextValue "C51" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D51" '7.52'
Set-TextValue "E51" '  -1.02%  '
